# Applies the "actualizacion de vistas de formularios en general" commit to
# Hoja1 (sheet1): updates the serial-number text in B233 (was "NA", now a
# real serial "SERT0002"), refreshes several "Cantidad" (quantity) values
# in column D, and updates the active view/selection to match where the
# user had scrolled/selected when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Real data corrections -------------------------------------------------

# B233 held the placeholder "NA" (no serial captured yet); it now holds the
# actual captured serial number.
$ws.Range("B233").Value = "SERT0002"

# Column D ("Cantidad") quantity corrections for the rows the user recounted.
$ws.Range("D8").Value = 20
$ws.Range("D9").Value = 30
$ws.Range("D15").Value = 20
$ws.Range("D17").Value = 20
$ws.Range("D18").Value = 20
$ws.Range("D19").Value = 30
$ws.Range("D20").Value = 15
$ws.Range("D22").Value = 25
$ws.Range("D25").Value = 20
$ws.Range("D35").Value = 15
$ws.Range("D43").Value = 20
$ws.Range("D45").Value = 50
$ws.Range("D46").Value = 40
$ws.Range("D47").Value = 25
$ws.Range("D49").Value = 25
$ws.Range("D52").Value = 20
$ws.Range("D58").Value = 20
$ws.Range("D62").Value = 30
$ws.Range("D66").Value = 15
$ws.Range("D77").Value = 15
$ws.Range("D84").Value = 20
$ws.Range("D85").Value = 30
$ws.Range("D88").Value = 20
$ws.Range("D89").Value = 20
$ws.Range("D90").Value = 15
$ws.Range("D96").Value = 30
$ws.Range("D103").Value = 25
$ws.Range("D110").Value = 25
$ws.Range("D115").Value = 25
$ws.Range("D126").Value = 20
$ws.Range("D138").Value = 30
$ws.Range("D154").Value = 25
$ws.Range("D155").Value = 30
$ws.Range("D157").Value = 15
$ws.Range("D164").Value = 15
$ws.Range("D165").Value = 20
$ws.Range("D166").Value = 50
$ws.Range("D167").Value = 50
$ws.Range("D174").Value = 20
$ws.Range("D184").Value = 30
$ws.Range("D189").Value = 25
$ws.Range("D190").Value = 25
$ws.Range("D192").Value = 25
$ws.Range("D199").Value = 30
$ws.Range("D200").Value = 20
$ws.Range("D201").Value = 25
$ws.Range("D202").Value = 50
$ws.Range("D212").Value = 40
$ws.Range("D213").Value = 30
$ws.Range("D219").Value = 35
$ws.Range("D221").Value = 35
$ws.Range("D222").Value = 35
$ws.Range("D223").Value = 35
$ws.Range("D224").Value = 20
$ws.Range("D226").Value = 20
$ws.Range("D229").Value = 30
$ws.Range("D232").Value = 30
$ws.Range("D241").Value = 20
$ws.Range("D242").Value = 30
$ws.Range("D248").Value = 30
$ws.Range("D250").Value = 20
$ws.Range("D251").Value = 30
$ws.Range("D252").Value = 20
$ws.Range("D253").Value = 15
$ws.Range("D255").Value = 35
$ws.Range("D258").Value = 15
$ws.Range("D262").Value = 25
$ws.Range("D271").Value = 25
$ws.Range("D272").Value = 30
$ws.Range("D278").Value = 30
$ws.Range("D279").Value = 20
$ws.Range("D280").Value = 30
$ws.Range("D281").Value = 20
$ws.Range("D282").Value = 30
$ws.Range("D283").Value = 30
$ws.Range("D285").Value = 25
$ws.Range("D288").Value = 25

# --- View / selection state --------------------------------------------

# The user scrolled further down the sheet and left the cursor on B234
# before saving.
$ws.Activate()
$ws.Range("A219").Select()
$excel.ActiveWindow.ScrollRow = 219
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B234").Select()
